$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns (E,F) in front of the existing "image" column,
#    which pushes it from E to G. This mirrors the sharedStrings/sheet1 diff
#    that introduces "FOB Price" / "CIF Price" columns ahead of the image
#    column.
# ---------------------------------------------------------------------------
$ws.Range("E1:F1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2. New header cells E1/F1 ("FOB Price" / "CIF Price"): start from the
#    "image" header's format (now in G1) which already has the yellow fill,
#    the header font and the thin left/right border, then switch the
#    alignment to centered like the "No" header.
# ---------------------------------------------------------------------------
$ws.Range("G1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E1").Value = "FOB Price"
$ws.Range("F1").Value = "CIF Price"
$ws.Range("E1:F1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("E1:F1").VerticalAlignment = -4160     # xlTop

# ---------------------------------------------------------------------------
# 3. New data cells E2/F2 (1000 / 2000): start from D2's number format
#    (thin border on all sides + the Chinese body font), then right align.
# ---------------------------------------------------------------------------
$ws.Range("D2").Copy()
$ws.Range("E2:F2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E2").Value = 1000
$ws.Range("F2").Value = 2000
$ws.Range("E2:F2").HorizontalAlignment = -4152   # xlRight
$ws.Range("E2:F2").VerticalAlignment = -4160     # xlTop

# ---------------------------------------------------------------------------
# 4. The header row's original columns (A1:D1) drop their bottom border so
#    the header visually merges with the new boxed columns.
# ---------------------------------------------------------------------------
$ws.Range("A1:D1").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none

# ---------------------------------------------------------------------------
# 5. The "insert image" cell (now G2) gains a full thin border on every side
#    and loses its left-aligned horizontal alignment (vertical=top stays).
# ---------------------------------------------------------------------------
$ws.Range("G2").Borders.LineStyle = 1
$ws.Range("G2").Borders.Weight = 2
$ws.Range("G2").HorizontalAlignment = 1          # xlGeneral
$ws.Range("G2").VerticalAlignment = -4160        # xlTop

# ---------------------------------------------------------------------------
# 6. Column widths: give the new FOB/CIF columns the same width as the
#    HSCode column and resize the rest to their new target widths.
# ---------------------------------------------------------------------------
$ws.Range("B1").ColumnWidth = 14.33
$ws.Range("C1").ColumnWidth = 20.33
$ws.Range("E1").ColumnWidth = 21.11
$ws.Range("F1").ColumnWidth = 21.11
$ws.Range("G1").ColumnWidth = 10.89

# ---------------------------------------------------------------------------
# 7. Selection cosmetics to match the saved workbook view.
# ---------------------------------------------------------------------------
$null = $ws.Range("E5").Select()
